$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates per the commit diff (Last/Previous figures + report dates refreshed).
# Numeric-looking values are prefixed with a literal apostrophe so Excel stores them as
# text (matching the workbook's existing inlineStr convention) instead of auto-converting
# them to numbers.

# Row 2
$ws.Range("C2").Value = "'7.25"
$ws.Range("D2").Value = "'7.24"

# Row 3
$ws.Range("C3").Value = "'3425"
$ws.Range("D3").Value = "'3439"

# Row 53
$ws.Range("C53").Value = "'0.2"
$ws.Range("D53").Value = "'0.1"
$ws.Range("H53").Value = "Oct/24"

# Row 60
$ws.Range("C60").Value = "'0"
$ws.Range("D60").Value = "'-0.1"
$ws.Range("H60").Value = "Oct/24"

# Row 61
$ws.Range("C61").Value = "'-0.1"
$ws.Range("H61").Value = "Oct/24"

# Row 62
$ws.Range("C62").Value = "'-0.3"
$ws.Range("D62").Value = "'-0.1"
$ws.Range("H62").Value = "Oct/24"

# Row 65
$ws.Range("D65").Value = "'1.86"

# Row 69
$ws.Range("C69").Value = "'500"
$ws.Range("D69").Value = "'1590"
$ws.Range("H69").Value = "Oct/24"

# Row 72
$ws.Range("C72").Value = "'14000"
$ws.Range("D72").Value = "'37600"
$ws.Range("H72").Value = "Oct/24"

# Row 73
$ws.Range("C73").Value = "'7.8"
$ws.Range("D73").Value = "'8.1"
$ws.Range("H73").Value = "Oct/24"

# Row 79
$ws.Range("C79").Value = "'233"
$ws.Range("D79").Value = "'126"

# Row 90
$ws.Range("C90").Value = "'641"
$ws.Range("D90").Value = "'580"
$ws.Range("H90").Value = "Sep/24"

# Row 93
$ws.Range("C93").Value = "'4215"
$ws.Range("D93").Value = "'4250"
$ws.Range("H93").Value = "Jul/24"

# Row 113
$ws.Range("C113").Value = "'1258"
$ws.Range("D113").Value = "'12113"
$ws.Range("H113").Value = "Sep/24"

# Row 114
$ws.Range("C114").Value = "'2788"
$ws.Range("D114").Value = "'18435"
$ws.Range("H114").Value = "Sep/24"

# Row 130
$ws.Range("C130").Value = "'3050000"
$ws.Range("D130").Value = "'2809000"
$ws.Range("H130").Value = "Oct/24"

# Row 131
$ws.Range("C131").Value = "'149"
$ws.Range("H131").Value = "Sep/24"
